# "updation of custom order"
# The QA test-data sheet "Sheet1" keeps a single tracked custom test order in row 2:
#   column A = OrderDate, column L = OverageID
# This commit refreshes that tracked order to a new date / overage id.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column A (OrderDate) and column L (OverageID) both hold plain text in this
# sheet (e.g. "12-01-2021", "58327790"). A straight .Value assignment would
# let Excel auto-detect these as a date / number, so a leading apostrophe is
# used to force them to stay literal text, exactly like typing them in by hand.
$ws.Range("A2").Formula = "'01-03-2022"
$ws.Range("L2").Formula = "'58532883"
